$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 410
$ws1.Range("F7").Value = 230
$ws1.Range("F8").Value = 1132
$ws1.Range("F9").Value = 312
$ws1.Range("F17").Value = 165
$ws1.Range("F18").Value = 2873
$ws1.Range("F26").Value = 5216
$ws1.Range("F28").Value = 975
$ws1.Range("F31").Value = 285

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1080
$ws2.Range("F5").Value = 1080
$ws2.Range("F27").Value = 3871
$ws2.Range("F32").Value = 42

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1773
$ws3.Range("F6").Value = 1015
$ws3.Range("F9").Value = 1277

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1773
$ws4.Range("F7").Value = 1015
$ws4.Range("F8").Value = 1277
$ws4.Range("F12").Value = 410
$ws4.Range("F14").Value = 230
$ws4.Range("F15").Value = 1132
$ws4.Range("F16").Value = 312
$ws4.Range("F19").Value = 1080
$ws4.Range("F23").Value = 165
$ws4.Range("F24").Value = 2873
$ws4.Range("F30").Value = 5216
$ws4.Range("F32").Value = 975
$ws4.Range("F38").Value = 285
$ws4.Range("F48").Value = 42
